# LeadZen_template.xlsx edit script
# Renames the "Financial Advisor" sheet to "Leads", rewrites the Index
# sheet's table-of-contents / details layout, and updates the Leads
# sheet's header text to match the new "Leads" export terminology.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Copy formats from their eventual-donor cells FIRST, while every
#    donor still has its original (pre-edit) formatting. This lets the
#    workbook's style table dedupe naturally instead of minting new
#    xf records for styles that already exist.
# ---------------------------------------------------------------------
$idx = $wb.Worksheets.Item("Index")

# G6:L6 becomes the new merged "Other Details" banner -- same visual
# style as the big left-aligned banners that used to live at G4:L4 and
# G14:L14.
$idx.Range("G4").Copy()
$idx.Range("G6:L6").PasteSpecial(-4122)

# K16 keeps the small blank "input box" style that used to sit at K24.
$idx.Range("K24").Copy()
$idx.Range("K16").PasteSpecial(-4122)

# G18 becomes the "View My Export History" link -- re-use the italic
# link style that used to live at G11 / G26.
$idx.Range("G26").Copy()
$idx.Range("G18").PasteSpecial(-4122)

# B22 / B23 become the footnote + copyright lines -- re-use the
# footnote/copyright styles that used to live at B29 / B30.
$idx.Range("B29").Copy()
$idx.Range("B22").PasteSpecial(-4122)
$idx.Range("B30").Copy()
$idx.Range("B23").PasteSpecial(-4122)

# B7 becomes a plain text line like B6 / B8.
$idx.Range("B6").Copy()
$idx.Range("B7:C7").PasteSpecial(-4122)

$idx.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Remove the "Search Details" banner box (G4:L4) entirely -- that
#    whole sub-section goes away with the Financial-Advisor search UI.
# ---------------------------------------------------------------------
$idx.Range("G4:L4").UnMerge()
$idx.Range("G4:L4").ClearContents()
$idx.Range("G4:L4").ClearFormats()

# Remove the old "Sort Applied / Date Added" and "Filter Applied / Feed
# Name Is Funded" detail rows and the search-results link underneath.
$idx.Range("J6").ClearContents()
$idx.Range("J6").ClearFormats()
$idx.Range("J8").ClearContents()
$idx.Range("J8").ClearFormats()
$idx.Range("J9").ClearContents()
$idx.Range("J9").ClearFormats()
$idx.Range("G11").ClearContents()
$idx.Range("G11").ClearFormats()

# The "Other Details" banner that used to live at G14:L14 moves up to
# G6:L6 (already formatted above), so clear out the old banner cells.
$idx.Range("G14:L14").UnMerge()
$idx.Range("G14:L14").ClearContents()
$idx.Range("G14:L14").ClearFormats()

# ---------------------------------------------------------------------
# 3. Re-title the table of contents entries.
# ---------------------------------------------------------------------
$idx.Range("B6").Value = "1 Leads"
$idx.Range("B7").Value = "2 Disclaimer"
$idx.Range("B7:C7").Merge()
$idx.Range("B8").ClearContents()

# ---------------------------------------------------------------------
# 4. Rebuild the "Other Details" panel, shifted two rows higher than
#    it used to be (it now starts right under the new G6:L6 banner).
# ---------------------------------------------------------------------
$idx.Range("G6:L6").Value = "Other Details"

$idx.Range("G8").Value = "Credits Utilised for this Export"

$idx.Range("G10").Value = "Number of Results"

$idx.Range("G12").Value = "Export Requested By"

$idx.Range("G14").Value = "Export Request Date "

$idx.Range("G16").Value = "Export ID"

$idx.Range("G18").Value = "View My Export History"

# ---------------------------------------------------------------------
# 5. Clear out everything below the old "Other Details" panel -- the
#    rows that used to hold G20/G22/G24/K24/G26 and the B29/B30
#    footnote + copyright (those got copied up to B22/B23 already).
# ---------------------------------------------------------------------
$idx.Range("G20").ClearContents()
$idx.Range("G20").ClearFormats()
$idx.Range("G22").ClearContents()
$idx.Range("G22").ClearFormats()
$idx.Range("G24").ClearContents()
$idx.Range("G24").ClearFormats()
$idx.Range("K24").ClearContents()
$idx.Range("K24").ClearFormats()
$idx.Range("G26").ClearContents()
$idx.Range("G26").ClearFormats()

$idx.Range("B22").Value = "** For enquiries concerning this export request, please contact support@leadzen.ai with your Export Id XXXXX"
$idx.Range("B23").Value = "Copyright © 2021, Leadzen.Ai  All rights reserved."
$idx.Range("B29").ClearContents()
$idx.Range("B29").ClearFormats()
$idx.Range("B30").ClearContents()
$idx.Range("B30").ClearFormats()

# Updated used range / cursor position.
$idx.Range("K18").Select() | Out-Null

# ---------------------------------------------------------------------
# 6. Rename the "Financial Advisor" sheet to "Leads" and update its
#    own header text (was a one-page "Financial Advisor" section
#    cover, now a "Leads" section cover).
# ---------------------------------------------------------------------
$leads = $wb.Worksheets.Item("Financial Advisor")
$leads.Name = "Leads"

$leads.Range("D1").Value = "Leads"
$leads.Range("D2").Value = "Details about the companies in your export set."
$leads.Range("D3").Value = "For internal use only"
$leads.Range("D4").Value = "Copyright © 2021, Leadzen.Ai  All rights reserved."

$leads.Range("C11").Select() | Out-Null

Write-Output "Index + Leads sheets updated"
